$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3: the "Files_to_upload" (D) paths now point to the
# new Testing_Env subfolder.
$ws.Range("D2").Value = "\Testdata\Templates\ImportPublications\Testing_Env\ICER RRMM - Success Case Sheet.xlsx"
$ws.Range("D3").Value = "\Testdata\Templates\ImportPublications\Testing_Env\mCRPC - Failure Case Sheet.xlsx"

# Add new regression scenario rows (LIVEHTA-67, 290, 517)
# Column A is populated first for all new rows (pop3..pop6)...
$ws.Range("A4").Value = "pop3"
$ws.Range("A5").Value = "pop4"
$ws.Range("A6").Value = "pop5"
$ws.Range("A7").Value = "pop6"

# ...then each row is completed, writing the Files_to_upload path (D) before
# the Expected_File_names (C) value, matching the author's editing order.
$ws.Range("B4").Value = "ICER RRMM 2022 report - ICER - Ovid search - 4/11/2022"
$ws.Range("D4").Value = "\Testdata\Templates\ImportPublications\Testing_Env\ICER RRMM - Header Mismatch.xlsx"
$ws.Range("C4").Value = "ICER RRMM - Header Mismatch.xlsx"

$ws.Range("B5").Value = "ICER RRMM 2022 report - ICER - Ovid search - 4/11/2022"
$ws.Range("D5").Value = "\Testdata\Templates\ImportPublications\Testing_Env\ICER RRMM - Letters in Publication Identifier.xlsx"
$ws.Range("C5").Value = "ICER RRMM - Letters in Publication Identifier.xlsx"

$ws.Range("B6").Value = "ICER RRMM 2022 report - ICER - Ovid search - 4/11/2022"
$ws.Range("D6").Value = "\Testdata\Templates\ImportPublications\Testing_Env\ICER RRMM - Empty value in Publication Identifier.xlsx"
$ws.Range("C6").Value = "ICER RRMM - Empty value in Publication Identifier.xlsx"

$ws.Range("B7").Value = "ICER RRMM 2022 report - ICER - Ovid search - 4/11/2022"
$ws.Range("D7").Value = "\Testdata\Templates\ImportPublications\Testing_Env\ICER RRMM - Duplicate value in FA18 column.xlsx"
$ws.Range("C7").Value = "ICER RRMM - Duplicate value in FA18 column.xlsx"

# Resize columns B and C to fit the new, wider content
$ws.Columns.Item(2).ColumnWidth = 48.21875
$ws.Columns.Item(3).ColumnWidth = 46

# Restore the selection to match the author's final cursor position
$ws.Range("C16").Select()
